# Scheduled-runner update: refresh computed profit-tracking columns (H-N)
# on a handful of rows across each class sheet. Values below reproduce the
# recalculated figures exactly (including clearing/adding the M/N "deficit"
# cell on rows where the recompute crossed the zero threshold).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 125
$ws.Range("I2").Value = 125
$ws.Range("K2").Value = 125
$ws.Range("M2").Value = -12

# Row 9
$ws.Range("H9").Value = 2250
$ws.Range("I9").Value = 2150.25
$ws.Range("K9").Value = 2150.25
$ws.Range("M9").Value = -1981.25

# Row 17
$ws.Range("H17").Value = 2190.3333
$ws.Range("J17").Value = 2190.3333
$ws.Range("L17").Value = 6570.999899999999
$ws.Range("N17").Value = -6906.999899999999

# Row 40
$ws.Range("H40").Value = 8708.25
$ws.Range("J40").Value = 10356.429
$ws.Range("L40").Value = 10356.429
$ws.Range("N40").Value = -10706.429

# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").Value = $null

# Row 54
$ws.Range("H54").Value = 9000
$ws.Range("I54").Value = 9000
$ws.Range("K54").Value = 9000
$ws.Range("M54").Value = -8514

# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = $null

# Row 76
$ws.Range("H76").Value = 83338970
$ws.Range("I76").Value = 4486.25
$ws.Range("K76").Value = 4486.25
$ws.Range("M76").Value = -4171.25

# Row 79
$ws.Range("H79").Value = 83338970
$ws.Range("I79").Value = 4486.25
$ws.Range("K79").Value = 4486.25
$ws.Range("M79").Value = -3394.25

# Row 86
$ws.Range("H86").Value = 5852008.5
$ws.Range("I86").Value = 2251
$ws.Range("J86").Value = 7523368
$ws.Range("K86").Value = 2251
$ws.Range("L86").Value = 7523368
$ws.Range("M86").Value = -1128
$ws.Range("N86").Value = -7525614

# Row 89
$ws.Range("H89").Value = 5852008.5
$ws.Range("I89").Value = 2251
$ws.Range("J89").Value = 7523368
$ws.Range("K89").Value = 11255
$ws.Range("L89").Value = 37616840
$ws.Range("M89").Value = -5639
$ws.Range("N89").Value = -37628072

# Row 107
$ws.Range("H107").Value = 31794.344
$ws.Range("I107").Value = 35045.586
$ws.Range("K107").Value = 35045.586
$ws.Range("M107").Value = -33125.586

# Row 112
$ws.Range("H112").Value = 2480.4333
$ws.Range("J112").Value = 2480.4333
$ws.Range("L112").Value = 7441.2999
$ws.Range("N112").Value = -9657.2999

# Row 125
$ws.Range("H125").Value = 22226246
$ws.Range("J125").Value = 27782252
$ws.Range("L125").Value = 250040268
$ws.Range("N125").Value = -250045188

# Row 135
$ws.Range("H135").Value = 716651
$ws.Range("I135").Value = 1001830.4
$ws.Range("J135").Value = 3702.5
$ws.Range("K135").Value = 9016473.6
$ws.Range("L135").Value = 33322.5
$ws.Range("M135").Value = -9013938.6
$ws.Range("N135").Value = -38392.5

# Row 137
$ws.Range("H137").Value = 459626.22
$ws.Range("I137").Value = 481179.84
$ws.Range("J137").Value = 7000
$ws.Range("K137").Value = 1443539.52
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = -1440989.52
$ws.Range("N137").Value = -26100

# Row 141
$ws.Range("H141").Value = 4527
$ws.Range("I141").Value = 3716.8
$ws.Range("K141").Value = 11150.4
$ws.Range("M141").Value = -5970.400000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 5187.25
$ws.Range("I6").Value = 5187.25
$ws.Range("K6").Value = 5187.25
$ws.Range("M6").Value = -5014.25

# Row 32
$ws.Range("H32").Value = 3599.53
$ws.Range("I32").Value = 3257.2405
$ws.Range("K32").Value = 3257.2405
$ws.Range("M32").Value = -2970.2405

# Row 61
$ws.Range("H61").Value = 3759
$ws.Range("I61").Value = 1971.8334
$ws.Range("J61").Value = 7333.3335
$ws.Range("K61").Value = 1971.8334
$ws.Range("L61").Value = 7333.3335
$ws.Range("M61").Value = -1759.8334
$ws.Range("N61").Value = -7757.3335

# Row 102
$ws.Range("H102").Value = 2130.5
$ws.Range("I102").Value = 2141.3635
$ws.Range("J102").Value = 2011
$ws.Range("K102").Value = 2141.3635
$ws.Range("L102").Value = 2011
$ws.Range("M102").Value = -519.3634999999999
$ws.Range("N102").Value = -5255

# Row 122
$ws.Range("H122").Value = 5465.522
$ws.Range("I122").Value = 6505
$ws.Range("K122").Value = 19515
$ws.Range("M122").Value = -17065

# Row 136
$ws.Range("H136").Value = 3759
$ws.Range("I136").Value = 1971.8334
$ws.Range("J136").Value = 7333.3335
$ws.Range("K136").Value = 5915.5002
$ws.Range("L136").Value = 22000.0005
$ws.Range("M136").Value = -3365.5002
$ws.Range("N136").Value = -27100.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1548135.1
$ws.Range("I86").Value = 2127186
$ws.Range("J86").Value = 3999.6667
$ws.Range("K86").Value = 2127186
$ws.Range("L86").Value = 3999.6667
$ws.Range("M86").Value = -2126063
$ws.Range("N86").Value = -6245.6667

# Row 89
$ws.Range("H89").Value = 1548135.1
$ws.Range("I89").Value = 2127186
$ws.Range("J89").Value = 3999.6667
$ws.Range("K89").Value = 10635930
$ws.Range("L89").Value = 19998.3335
$ws.Range("M89").Value = -10630314
$ws.Range("N89").Value = -31230.3335

# Row 99
$ws.Range("H99").Value = 5067.857
$ws.Range("I99").Value = 7116
$ws.Range("J99").Value = 2337
$ws.Range("K99").Value = 7116
$ws.Range("L99").Value = 2337
$ws.Range("M99").Value = -5618
$ws.Range("N99").Value = -5333

# Row 107
$ws.Range("H107").Value = 418202.84
$ws.Range("I107").Value = 1091.5
$ws.Range("K107").Value = 1091.5
$ws.Range("M107").Value = 828.5

# Row 134
$ws.Range("H134").Value = 108818.9
$ws.Range("I134").Value = 9741.714
$ws.Range("J134").Value = 339999
$ws.Range("K134").Value = 29225.142
$ws.Range("L134").Value = 1019997
$ws.Range("M134").Value = -26690.142
$ws.Range("N134").Value = -1025067

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 4348.5
$ws.Range("I99").Value = 4148
$ws.Range("J99").Value = 4749.5
$ws.Range("K99").Value = 4148
$ws.Range("L99").Value = 4749.5
$ws.Range("M99").Value = -2650
$ws.Range("N99").Value = -7745.5

# Row 126
$ws.Range("H126").Value = 4348.5
$ws.Range("I126").Value = 4148
$ws.Range("J126").Value = 4749.5
$ws.Range("K126").Value = 12444
$ws.Range("L126").Value = 14248.5
$ws.Range("M126").Value = -9974
$ws.Range("N126").Value = -19188.5

# Row 141
$ws.Range("H141").Value = 525218.8
$ws.Range("J141").Value = 626449.5
$ws.Range("L141").Value = 626449.5
$ws.Range("N141").Value = -636809.5

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1951001.2
$ws.Range("J113").Value = 2045.3846
$ws.Range("L113").Value = 6136.1538
$ws.Range("N113").Value = -10476.1538

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 5025000

# Row 55
$ws.Range("H55").Value = 18751.666
$ws.Range("J55").Value = 18125
$ws.Range("L55").Value = 18125
$ws.Range("N55").Value = -18779

# Row 80
$ws.Range("H80").Value = 2004133.2
$ws.Range("I80").Value = 1433648.9
$ws.Range("J80").Value = 3335263.2
$ws.Range("K80").Value = 1433648.9
$ws.Range("L80").Value = 3335263.2
$ws.Range("M80").Value = -1432650.9
$ws.Range("N80").Value = -3337259.2

# Row 83
$ws.Range("H83").Value = 2004133.2
$ws.Range("I83").Value = 1433648.9
$ws.Range("J83").Value = 3335263.2
$ws.Range("K83").Value = 7168244.5
$ws.Range("L83").Value = 16676316
$ws.Range("M83").Value = -7163252.5
$ws.Range("N83").Value = -16686300

# Row 97
$ws.Range("H97").Value = 816.4211
$ws.Range("I97").Value = 816.4211
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 816.4211
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -320.4211
$ws.Range("N97").Value = $null

# Row 102
$ws.Range("H102").Value = 1954.9
$ws.Range("I102").Value = 969.8461
$ws.Range("J102").Value = 3784.2856
$ws.Range("K102").Value = 969.8461
$ws.Range("L102").Value = 3784.2856
$ws.Range("M102").Value = 652.1539
$ws.Range("N102").Value = -7028.2856

# Row 122
$ws.Range("H122").Value = 3767.6667
$ws.Range("I122").Value = 1844.5714
$ws.Range("K122").Value = 5533.7142
$ws.Range("M122").Value = -3083.7142

# Row 132
$ws.Range("H132").Value = 62888.832
$ws.Range("I132").Value = 7499.1665
$ws.Range("J132").Value = 173668.17
$ws.Range("K132").Value = 22497.4995
$ws.Range("L132").Value = 521004.51
$ws.Range("M132").Value = -19967.4995
$ws.Range("N132").Value = -526064.51

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 83335256
$ws.Range("I93").Value = 111112744
$ws.Range("J93").Value = 2766.3333
$ws.Range("K93").Value = 111112744
$ws.Range("L93").Value = 2766.3333
$ws.Range("M93").Value = -111111496
$ws.Range("N93").Value = -5262.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1306.84
$ws.Range("I81").Value = 1359.7391
$ws.Range("J81").Value = 698.5
$ws.Range("K81").Value = 2719.4782
$ws.Range("L81").Value = 1397
$ws.Range("M81").Value = -1658.4782
$ws.Range("N81").Value = -3519

# Row 84
$ws.Range("H84").Value = 1306.84
$ws.Range("I84").Value = 1359.7391
$ws.Range("J84").Value = 698.5
$ws.Range("K84").Value = 13597.391
$ws.Range("L84").Value = 6985
$ws.Range("M84").Value = -8293.391
$ws.Range("N84").Value = -17593

# Row 132
$ws.Range("H132").Value = 29618.5
$ws.Range("I132").Value = 2348.6785
$ws.Range("K132").Value = 7046.0355
$ws.Range("M132").Value = -4516.0355
